# Update recomputed NATMI ligand/receptor expression & specificity metrics
# (columns G:T, rows 2:9) with newly generated TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = New-Object 'object[,]' 8,14
$newData[0,0] = 0.7454756666666666
$newData[0,1] = 2.236427
$newData[0,2] = 0.004237455828433692
$newData[0,3] = 0.004251944035061194
$newData[0,4] = 1
$newData[0,5] = 0.3333333333333333
$newData[0,6] = 0.1375686666666667
$newData[0,7] = 0.412706
$newData[0,8] = 0.2896572731203081
$newData[0,9] = 0.2896572731203081
$newData[0,10] = 0.1025540934957778
$newData[0,11] = 0.922986841462
$newData[0,12] = 0.001227409900231859
$newData[0,13] = 0.001231606514655985
$newData[1,0] = 0.7454756666666666
$newData[1,1] = 2.236427
$newData[1,2] = 0.004237455828433692
$newData[1,3] = 0.004251944035061194
$newData[1,4] = 2
$newData[1,5] = 0.6666666666666666
$newData[1,6] = 0.3373673333333334
$newData[1,7] = 1.012102
$newData[1,8] = 0.7103427268796919
$newData[1,9] = 0.7103427268796919
$newData[1,10] = 0.2514991377282222
$newData[1,11] = 2.263492239554
$newData[1,12] = 0.003010045928201832
$newData[1,13] = 0.003020337520405209
$newData[2,0] = 102.35201
$newData[2,1] = 307.05603
$newData[2,2] = 0.5817924591230612
$newData[2,3] = 0.583781654929077
$newData[2,4] = 1
$newData[2,5] = 0.3333333333333333
$newData[2,6] = 0.1375686666666667
$newData[2,7] = 0.412706
$newData[2,8] = 0.2896572731203081
$newData[2,9] = 0.2896572731203081
$newData[2,10] = 14.08042954635333
$newData[2,11] = 126.72386591718
$newData[2,12] = 0.1685204172315443
$newData[2,13] = 0.1690966022644171
$newData[3,0] = 102.35201
$newData[3,1] = 307.05603
$newData[3,2] = 0.5817924591230612
$newData[3,3] = 0.583781654929077
$newData[3,4] = 2
$newData[3,5] = 0.6666666666666666
$newData[3,6] = 0.3373673333333334
$newData[3,7] = 1.012102
$newData[3,8] = 0.7103427268796919
$newData[3,9] = 0.7103427268796919
$newData[3,10] = 34.53022467500666
$newData[3,11] = 310.77202207506
$newData[3,12] = 0.413272041891517
$newData[3,13] = 0.4146850526646599
$newData[4,0] = 1.79836
$newData[4,1] = 3.59672
$newData[4,2] = 0.01022229350247785
$newData[4,3] = 0.006838162904394061
$newData[4,4] = 1
$newData[4,5] = 0.3333333333333333
$newData[4,6] = 0.1375686666666667
$newData[4,7] = 0.412706
$newData[4,8] = 0.2896572731203081
$newData[4,9] = 0.2896572731203081
$newData[4,10] = 0.2473979873866667
$newData[4,11] = 1.48438792432
$newData[4,12] = 0.002960961660963179
$newData[4,13] = 0.00198072362003923
$newData[5,0] = 1.79836
$newData[5,1] = 3.59672
$newData[5,2] = 0.01022229350247785
$newData[5,3] = 0.006838162904394061
$newData[5,4] = 2
$newData[5,5] = 0.6666666666666666
$newData[5,6] = 0.3373673333333334
$newData[5,7] = 1.012102
$newData[5,8] = 0.7103427268796919
$newData[5,9] = 0.7103427268796919
$newData[5,10] = 0.6067079175733333
$newData[5,11] = 3.64024750544
$newData[5,12] = 0.007261331841514676
$newData[5,13] = 0.004857439284354831
$newData[6,0] = 71.02944933333333
$newData[6,1] = 213.088348
$newData[6,2] = 0.4037477915460271
$newData[6,3] = 0.4051282381314676
$newData[6,4] = 1
$newData[6,5] = 0.3333333333333333
$newData[6,6] = 0.1375686666666667
$newData[6,7] = 0.412706
$newData[6,8] = 0.2896572731203081
$newData[6,9] = 0.2896572731203081
$newData[6,10] = 9.771426638854223
$newData[6,11] = 87.942839749688
$newData[6,12] = 0.1169484843275688
$newData[6,13] = 0.1173483407211958
$newData[7,0] = 71.02944933333333
$newData[7,1] = 213.088348
$newData[7,2] = 0.4037477915460271
$newData[7,3] = 0.4051282381314676
$newData[7,4] = 2
$newData[7,5] = 0.6666666666666666
$newData[7,6] = 0.3373673333333334
$newData[7,7] = 1.012102
$newData[7,8] = 0.7103427268796919
$newData[7,9] = 0.7103427268796919
$newData[7,10] = 23.96301590972178
$newData[7,11] = 215.667143187496
$newData[7,12] = 0.2867993072184583
$newData[7,13] = 0.2877798974102719

$ws.Range("G2:T9").Value = $newData
